$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.375
$ws.Range("C2").Value = 6

$ws.Range("B3").Value = 5.375
$ws.Range("C3").Value = 6

$ws.Range("B4").Value = 6.875
$ws.Range("C4").Value = 7

$ws.Range("B5").Value = 12.875
$ws.Range("C5").Value = 13

$ws.Range("B6").Value = 34.04166666666666
$ws.Range("C6").Value = 35

$ws.Range("B7").Value = 6.875
$ws.Range("C7").Value = 7
